$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.192.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.30%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.262.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.91%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'496.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.25%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'128.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.23%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.33%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.525"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.55%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.0951"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.86%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.82%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.336"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.41%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.50%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.662.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.30%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'22.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.90%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'54.164.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.13%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000130"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.75%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.261.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.11%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'10.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.51%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.43%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'302.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.32%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.42%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.21%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'60.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.19%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.19%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.149"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.81%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'170.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.37%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.70%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0691"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.57%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.08%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.22%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.16%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'17.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.17%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.67%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.944"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +8.85%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.49%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.51%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.373"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.34%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.31%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.74%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'125.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.47%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0494"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.24%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0891"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.545"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.67%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'241.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.63%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.372"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.16%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0204"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.04%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.35%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'16.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.69%  "
$ws.Range("E51").Style = "Normal"
